# trading_journal.xlsx update:
#  - fill in Exit/Result/P&L/Exit Reason/Status for rows 24 and 26 (trades closed)
#  - fill in Exit Reason for row 28
#  - fill in Trade Type/Currency for row 30
#  - add two brand-new trades in rows 31 and 32 (Castrol / Schneider Electric Infra)
#  - extend selection to A1:U32

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 24 - Castrol trade closed at a loss (SL Hit)
# ---------------------------------------------------------------------------
$ws.Range("G24").Value = 1907.8
$ws.Range("K24").Value = "Loss"
$ws.Range("L24").Formula = "=(G24-D24)*H24"
$ws.Range("N24").Value = "SL Hit"
$ws.Range("T24").Value = "Closed"

# ---------------------------------------------------------------------------
# Row 26 - trade closed at a profit (Take Profit Hit)
# ---------------------------------------------------------------------------
$ws.Range("G26").Value = 1857.1
$ws.Range("K26").Value = "Profit"
$ws.Range("L26").Formula = "=(G26-D26)*H26"
$ws.Range("N26").Value = "Trailing Stop Loss Hit"
$ws.Range("T26").Value = "Closed"

# ---------------------------------------------------------------------------
# Row 28 - add Exit Reason
# ---------------------------------------------------------------------------
$ws.Range("N28").Value = "SL Hit"

# ---------------------------------------------------------------------------
# Row 30 - add Trade Type / Currency
# ---------------------------------------------------------------------------
$ws.Range("Q30").Value = "Weekly Day"
$ws.Range("R30").Value = "INR"

# ---------------------------------------------------------------------------
# Row 31 - new trade: Castrol, Diamod Pattern
# ---------------------------------------------------------------------------
$ws.Range("A22").Copy()
$ws.Range("A31").PasteSpecial(-4122)
$ws.Range("A31").Value = 45868
$ws.Range("B31").Value = "Castrol"
$ws.Range("C31").Value = "Long"
$ws.Range("D31").Value = 220
$ws.Range("E31").Value = 207
$ws.Range("F31").Value = 272
$ws.Range("H31").Value = 231
$ws.Range("I31").Formula = "=(D31-E31)*H31"
$ws.Range("J31").Formula = "=(F31-D31)/(D31-E31)"
$ws.Range("M31").Value = "Diamod Pattern"
$ws.Range("Q31").Value = "Weekly Day"
$ws.Range("R31").Value = "INR"
$ws.Range("S31").Formula = "=L31"
$ws.Range("T31").Value = "Active"
$ws.Range("U31").Formula = "=H31*D31"

# ---------------------------------------------------------------------------
# Row 32 - new trade: Schneider Electric Infra, closed at a loss
# ---------------------------------------------------------------------------
$ws.Range("A22").Copy()
$ws.Range("A32").PasteSpecial(-4122)
$ws.Range("A32").Value = 45866
$ws.Range("B32").Value = "Schneider Electric Infra"
$ws.Range("C32").Value = "Long"
$ws.Range("D32").Value = 934
$ws.Range("E32").Value = 883.25
$ws.Range("F32").Value = 1100
$ws.Range("G32").Value = 883.25
$ws.Range("H32").Value = 63
$ws.Range("I32").Formula = "=(D32-E32)*H32"
$ws.Range("J32").Formula = "=(F32-D32)/(D32-E32)"
$ws.Range("K32").Value = "Loss"
$ws.Range("L32").Formula = "=(G32-D32)*H32"
$ws.Range("M32").Value = "All time Breakout "
$ws.Range("N32").Value = "SL Hit"
$ws.Range("Q32").Value = "Weekly Day"
$ws.Range("R32").Value = "INR"
$ws.Range("S32").Formula = "=L32"
$ws.Range("T32").Value = "Closed"
$ws.Range("U32").Formula = "=H32*D32"

# ---------------------------------------------------------------------------
# Selection / view
# ---------------------------------------------------------------------------
$ws.Range("A1:U32").Select()
